$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'MATEMATICA DISCRETA'
$ws.Range("B2").Value = 30
$ws.Range("A3").Value = 'METODOS MATEMATICOS I'
$ws.Range("B3").Value = 159
$ws.Range("A4").Value = 'PROGRAMACION'
$ws.Range("B4").Value = 145
$ws.Range("A5").Value = 'SEMINARIO DE SOLUCION DE PROBLEMAS DE PROGRAMACION'
$ws.Range("B5").Value = 123
$ws.Range("A6").Value = 'SEMINARIO DE SOLUCION DE PROBLEMAS DE METODOS MATEMATICOS I'
$ws.Range("B6").Value = 169
$ws.Range("A7").Value = 'ALGORITMIA'
$ws.Range("B7").Value = 36
$ws.Range("A8").Value = 'ESTRUCTURAS DE DATOS I'
$ws.Range("B8").Value = 0
$ws.Range("A9").Value = 'METODOS MATEMATICOS II'
$ws.Range("B9").Value = 134
$ws.Range("A10").Value = 'SEMINARIO DE SOLUCION DE PROBLEMAS DE ALGORITMIA'
$ws.Range("B10").Value = 117
$ws.Range("A11").Value = 'SEMINARIO DE SOLUCION DE PROBLEMAS DE ESTRUCTURAS DE DATOS I'
$ws.Range("B11").Value = 28
$ws.Range("A12").Value = 'SEMINARIO DE SOLUCION DE PROBLEMAS DE METODOS MATEMATICOS II'
$ws.Range("B12").Value = 123
$ws.Range("A13").Value = 'ADMINISTRACION DE REDES'
$ws.Range("B13").Value = 61
$ws.Range("A14").Value = 'ESTADISTICA Y PROCESOS ESTOCASTICOS'
$ws.Range("B14").Value = 48
$ws.Range("A15").Value = 'ESTRUCTURAS DE DATOS II'
$ws.Range("B15").Value = 0
$ws.Range("A16").Value = 'SEMINARIO DE SOLUCION DE PROBLEMAS DE ESTRUCTURAS DE DATOS II'
$ws.Range("B16").Value = 60
$ws.Range("A17").Value = 'TEORIA DE LA COMPUTACION'
$ws.Range("B17").Value = 32
$ws.Range("A18").Value = 'ADMINISTRACION DE SERVIDORES'
$ws.Range("B18").Value = 90
$ws.Range("A19").Value = 'BASES DE DATOS'
$ws.Range("B19").Value = 57
$ws.Range("A20").Value = 'HIPERMEDIA'
$ws.Range("B20").Value = 80
$ws.Range("A21").Value = 'INGENIERIA DE SOFTWARE I'
$ws.Range("B21").Value = 20
$ws.Range("A22").Value = 'SEMINARIO DE SOLUCION DE PROBLEMAS DE BASES DE DATOS'
$ws.Range("B22").Value = 128
$ws.Range("A23").Value = 'CONTROL DE PROYECTOS'
$ws.Range("B23").Value = 89
$ws.Range("A24").Value = 'INGENIERIA DE SOFTWARE II'
$ws.Range("B24").Value = 98
$ws.Range("A25").Value = 'PROGRAMACION PARA INTERNET'
$ws.Range("B25").Value = 196
$ws.Range("A26").Value = 'SEGURIDAD DE LA INFORMACION'
$ws.Range("B26").Value = 136
$ws.Range("A27").Value = 'SEMINARIO DE SOLUCION DE PROBLEMAS DE INGENIERIA DE SOFTWARE I'
$ws.Range("B27").Value = 148
$ws.Range("A28").Value = 'SEMINARIO DE SOLUCION DE PROBLEMAS DE USO, ADAPTACION, EXPLOTACION DE SISTEMAS OPERATIVOS'
$ws.Range("B28").Value = 116
$ws.Range("A29").Value = 'USO, ADAPTACION Y EXPLOTACION DE SISTEMAS OPERATIVOS'
$ws.Range("B29").Value = 105
$ws.Range("A30").Value = 'ADMINISTRACION DE BASES DE DATOS'
$ws.Range("B30").Value = 96
$ws.Range("A31").Value = 'ALMACENES DE DATOS (DATA WAREHOUSE)'
$ws.Range("B31").Value = 88
$ws.Range("A32").Value = 'MINERIA DE DATOS'
$ws.Range("B32").Value = 117
$ws.Range("A33").Value = 'CLASIFICACION INTELIGENTE DE DATOS'
$ws.Range("B33").Value = 99
$ws.Range("A34").Value = 'SEMINARIO DE SOLUCION DE PROBLEMAS DE SISTEMAS BASADOS EN CONOCIMIENTO'
$ws.Range("B34").Value = 93
$ws.Range("A35").Value = 'SISTEMAS BASADOS EN CONOCIMIENTO'
$ws.Range("B35").Value = 93

